# Update randomized test result statistics on the "stats" worksheet.
# The merged framework now re-checks all feasible paths for a recently
# merged component before the next nomination phase, which changes the
# recorded gmin-related statistics produced by the test run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

# --- Block 1 (rows 2-6) ---
$ws.Range("D2").Value = 0.0001074313186109066
$ws.Range("E2").Value = 0.03318778006359935
$ws.Range("G2").Value = 0.002988782245665789
$ws.Range("H2").Value = 0.005840327125042677
$ws.Range("I2").Value = 0.007988458499312401
$ws.Range("J2").Value = 0.01216139039024711
$ws.Range("K2").Value = 0.001206408720463514

$ws.Range("D3").Value = 0.002322259824723005
$ws.Range("E3").Value = 0.04118812456727028
$ws.Range("G3").Value = 0.00366108538582921
$ws.Range("H3").Value = 0.008668153546750546
$ws.Range("I3").Value = 0.008431130088865757
$ws.Range("J3").Value = 0.01543169608339667
$ws.Range("K3").Value = 0.001841188874095678

$ws.Range("D4").Value = 0.002208093646913767
$ws.Range("E4").Value = 0.0626468462869525
$ws.Range("G4").Value = 0.00326509540900588
$ws.Range("H4").Value = 0.0222645215690136
$ws.Range("I4").Value = 0.00873554265126586
$ws.Range("J4").Value = 0.02402074309065938
$ws.Range("K4").Value = 0.001227786298841238

$ws.Range("D5").Value = 0.0002701361663639545
$ws.Range("E5").Value = 0.07227285765111446
$ws.Range("G5").Value = 0.002981152851134539
$ws.Range("H5").Value = 0.01089310133829713
$ws.Range("I5").Value = 0.01989418035373092
$ws.Range("J5").Value = 0.03430342068895698
$ws.Range("K5").Value = 0.001179954502731562

$ws.Range("D6").Value = 0.003749554045498371
$ws.Range("E6").Value = 0.1677141105756164
$ws.Range("G6").Value = 0.006936208344995975
$ws.Range("H6").Value = 0.01874947315081954
$ws.Range("I6").Value = 0.1152381421998143
$ws.Range("J6").Value = 0.01773395063355565
$ws.Range("K6").Value = 0.002564535941928625

# --- Block 2 (rows 8-12), mirrors block 1 ---
$ws.Range("D8").Value = 0.0001074313186109066
$ws.Range("E8").Value = 0.03318778006359935
$ws.Range("G8").Value = 0.002988782245665789
$ws.Range("H8").Value = 0.005840327125042677
$ws.Range("I8").Value = 0.007988458499312401
$ws.Range("J8").Value = 0.01216139039024711
$ws.Range("K8").Value = 0.001206408720463514

$ws.Range("D9").Value = 0.002322259824723005
$ws.Range("E9").Value = 0.04118812456727028
$ws.Range("G9").Value = 0.00366108538582921
$ws.Range("H9").Value = 0.008668153546750546
$ws.Range("I9").Value = 0.008431130088865757
$ws.Range("J9").Value = 0.01543169608339667
$ws.Range("K9").Value = 0.001841188874095678

$ws.Range("D10").Value = 0.002208093646913767
$ws.Range("E10").Value = 0.0626468462869525
$ws.Range("G10").Value = 0.00326509540900588
$ws.Range("H10").Value = 0.0222645215690136
$ws.Range("I10").Value = 0.00873554265126586
$ws.Range("J10").Value = 0.02402074309065938
$ws.Range("K10").Value = 0.001227786298841238

$ws.Range("D11").Value = 0.0002701361663639545
$ws.Range("E11").Value = 0.07227285765111446
$ws.Range("G11").Value = 0.002981152851134539
$ws.Range("H11").Value = 0.01089310133829713
$ws.Range("I11").Value = 0.01989418035373092
$ws.Range("J11").Value = 0.03430342068895698
$ws.Range("K11").Value = 0.001179954502731562

$ws.Range("D12").Value = 0.003749554045498371
$ws.Range("E12").Value = 0.1677141105756164
$ws.Range("G12").Value = 0.006936208344995975
$ws.Range("H12").Value = 0.01874947315081954
$ws.Range("I12").Value = 0.1152381421998143
$ws.Range("J12").Value = 0.01773395063355565
$ws.Range("K12").Value = 0.002564535941928625
